$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 17: update F17 / G17 ---
$ws.Range("F17").Value = 40
$ws.Range("G17").Value = 40

# --- Row 50: fill in person percentages ---
$ws.Range("E50").Value = 10
$ws.Range("F50").Value = 10
$ws.Range("G50").Value = 70
$ws.Range("H50").Value = 10

# --- Row 56: new document row ---
$ws.Range("B56").Value = "reflexion_danielschmidt.doc"
$ws.Range("E56").Value = 0
$ws.Range("F56").Value = 100
$ws.Range("G56").Value = 0
$ws.Range("H56").Value = 0

# --- Sheet view: drop the scrolled top-left cell and move the selection ---
$ws.Range("F3:F11").Select()

$wb.Save()
